# Final-output pass: refresh column B (ORTHO TARGET), column C (PRODUCTION)
# and column E (Condition) for every data row with the newly computed
# values, per "completed the final output portion" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, B (ORTHO TARGET), C (PRODUCTION), E (Condition)
$rows = @(
    @(2, "Kinesin", "low", "low"),
    @(3, "Reducer", "s", "high"),
    @(4, "Isotope", "s", "low"),
    @(5, "Oxidizer", "s", "low"),
    @(6, "Eukaryote", "s", "low"),
    @(7, "Phosphorylation", "s", "high"),
    @(8, "Cytoplasm", "s", "low"),
    @(9, "Ribosomes", "s", "low"),
    @(10, "Glycolysis", "s", "low"),
    @(11, "Chlorophyll", "s", "low"),
    @(12, "Dynein", "s", "high"),
    @(13, "Isotonic", "s", "low"),
    @(14, "Nucleoid", "s", "high"),
    @(15, "Tertiary", "s", "low"),
    @(16, "Exergonic", "s", "low"),
    @(17, "Ribonucleoside", "s", "low"),
    @(18, "Purine", "s", "high"),
    @(19, "Quaternary", "s", "high"),
    @(20, "Catalyst", "s", "high"),
    @(21, "Desmosomes", "s", "high"),
    @(22, "Amphipathic", "s", "high"),
    @(23, "Monosaccharides", "s", "low"),
    @(24, "Peroxisome", "s", "low"),
    @(25, "Chemiosmotic", "s", "low"),
    @(26, "Hypertonic", "s", "high"),
    @(27, "Microtubule", "s", "high"),
    @(28, "Centrioles", "s", "low"),
    @(29, "Vacuole", "s", "high"),
    @(30, "Lysosome", "s", "high"),
    @(31, "Nucleotides", "s", "high")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]
    $ws.Cells.Item($rowNum, 3).Value = $r[2]
    $ws.Cells.Item($rowNum, 5).Value = $r[3]
}
